# Applies the commit "Update gh-pages to output generated at 456a3b4" to
# the workbook.  The scraper re-ran and:
#   - two events that had happened already (瑞昌·铜源剪影动漫游戏节 and the
#     "火只木南专场见面会" sub-event) dropped off the top of the list
#   - every remaining row shifted up by two positions
#   - the "想去人数" (F column) counters ticked up for most remaining events
# This affects worksheet 1 (展览) and worksheet 4 (全部类型); worksheets 2
# (演出) and 3 (本地生活) are untouched by the edit.

$wb = $excel.ActiveWorkbook

function Update-EventSheet {
    param([object]$ws, [object]$fValues)

    # Remove the two obsolete rows at the top of the data (rows 2 and 3);
    # this shifts all the remaining data rows up by two and keeps all of
    # their other column values (B-I) intact/correct.
    $ws.Rows.Item(2).Delete() | Out-Null
    $ws.Rows.Item(2).Delete() | Out-Null

    $rowCount = $fValues.Length

    for ($i = 0; $i -lt $rowCount; $i++) {
        $r = $i + 2
        # Column A holds the plain row index (row number - 1); row deletion
        # does not renumber literal values, so reset them to be sequential.
        $ws.Cells.Item($r, 1).Value = ($r - 1)
        # Column F (想去人数) needs refreshing to the newly scraped counts.
        $ws.Cells.Item($r, 6).Value = $fValues[$i]
    }
}

# Worksheet 1: 展览 (A1:I19 -> A1:I17, 16 data rows)
$ws1 = $wb.Worksheets.Item(1)
$f1 = @(4436, 842, 38, 135, 152, 33, 18, 137, 614, 23, 184, 1201, 19, 2818, 441, 513)
Update-EventSheet $ws1 $f1

# Worksheet 4: 全部类型 (A1:I20 -> A1:I18, 17 data rows)
$ws4 = $wb.Worksheets.Item(4)
$f4 = @(4436, 842, 38, 135, 152, 33, 18, 137, 614, 23, 1, 184, 1201, 19, 2818, 441, 513)
Update-EventSheet $ws4 $f4

Write-Output "Done"
